$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.747.61'
$ws.Range('E2').Value = '  +3.63%  '
$ws.Range('D3').Value = '2.269.71'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''303.76'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '''92.04'
$ws.Range('E6').Value = '  +4.40%  '
$ws.Range('E7').Value = '  +2.55%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.480'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('D10').Value = '''32.31'
$ws.Range('E10').Value = '  +4.01%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '''6.59'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '2.621.60'
$ws.Range('E15').Value = '  +2.46%  '
$ws.Range('D16').Value = '''14.19'
$ws.Range('E16').Value = '  +2.16%  '
$ws.Range('D17').Value = '2.283.53'
$ws.Range('E17').Value = '  +3.55%  '
$ws.Range('D18').Value = '''0.763'
$ws.Range('E18').Value = '  +3.35%  '
$ws.Range('D19').Value = '41.677.50'
$ws.Range('E19').Value = '  +3.80%  '
$ws.Range('D20').Value = '''12.56'
$ws.Range('E20').Value = '  +9.98%  '
$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').Value = '''5.91'
$ws.Range('E22').Value = '  +1.92%  '
$ws.Range('D23').Value = '''66.82'
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('D24').Value = '''240.03'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('E25').Value = '  +3.01%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  +5.04%  '
$ws.Range('D28').Value = '''24.04'
$ws.Range('E28').Value = '  +2.73%  '
$ws.Range('D29').Value = '''9.51'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('D31').Value = '''160.58'
$ws.Range('E31').Value = '  +2.36%  '
$ws.Range('D32').Value = '''34.13'
$ws.Range('E32').Value = '  +4.96%  '
$ws.Range('D33').Value = '''5.24'
$ws.Range('E33').Value = '  +4.81%  '
$ws.Range('D34').Value = '''1.00'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '''0.0743'
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').Value = '''16.85'
$ws.Range('E37').Value = '  +7.91%  '
$ws.Range('D38').Value = '''2.38'
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '''0.116'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.104'
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('D41').Value = '''1.79'
$ws.Range('E41').Value = '  +3.36%  '
$ws.Range('E42').Value = '  +2.95%  '
$ws.Range('D43').Value = '2.043.71'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('D44').Value = '''19.25'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = '''10.37'
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('D46').Value = '''0.0278'
$ws.Range('E46').Value = '  +2.39%  '
$ws.Range('E47').Value = '  +11.59%  '
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('D49').Value = '''1.53'
$ws.Range('E49').Value = '  +3.86%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '''1.16'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '''72.61'
$ws.Range('E51').Value = '  +6.16%  '
